$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Left table (columns A:H), data rows 3..10 ("negative" anchor) ----
$leftData = @(
    @("crude", 0.9411764705882353, 32, 32, 0, 1, $false, 2),
    @("fraud", 0.6944444444444444, 25, 25, 0, 1, $false, 11),
    @("crisis", 0.6643835616438356, 194, 194, 0, 1, $false, 98),
    @("panic", 0.2887596899224806, 149, 149, 0, 1, $false, 367),
    @("low", 0.2348993288590604, 35, 35, 0, 1, $false, 114),
    @("sc", 0.2275132275132275, 43, 43, 0, 1, $false, 146),
    @("stop", 0.1666666666666667, 42, 42, 0, 1, $false, 210),
    @("demand", 0.06720430107526881, 25, 26, 0.04, 0.96, $true, 347)
)

for ($i = 0; $i -lt $leftData.Length; $i++) {
    $r = 3 + $i
    $row = $leftData[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]   # A - name
    $ws.Cells.Item($r, 2).Value = $row[1]   # B - anchor score
    $ws.Cells.Item($r, 3).Value = $row[2]   # C - type occurences
    $ws.Cells.Item($r, 4).Value = $row[3]   # D - total occurences
    $ws.Cells.Item($r, 5).Value = $row[4]   # E - +%
    $ws.Cells.Item($r, 6).Value = $row[5]   # F - -%
    $ws.Cells.Item($r, 7).Value = $row[6]   # G - both
    $ws.Cells.Item($r, 8).Value = $row[7]   # H - normal
}

# Row 11 of the left table no longer has data (left table shrank from 9 to 8 words) - clear it.
$ws.Range("A11:H11").ClearContents()

# ---- Right table (columns J:Q), data rows 3..42 ("positive" anchor) ----
$rightData = @(
    @("happy", 0.9615384615384616, 25, 25, 1, 0, $false, 1),
    @("best", 0.9491525423728814, 56, 56, 1, 0, $false, 3),
    @("interesting", 0.9393939393939394, 31, 31, 1, 0, $false, 2),
    @("love", 0.9130434782608695, 42, 42, 1, 0, $false, 4),
    @("great", 0.8482142857142857, 95, 95, 1, 0, $false, 17),
    @("thanks", 0.8414634146341463, 69, 69, 1, 0, $false, 13),
    @("positive", 0.7931034482758621, 46, 46, 1, 0, $false, 12),
    @("thank", 0.7890625, 101, 101, 1, 0, $false, 27),
    @("free", 0.775, 93, 93, 1, 0, $false, 27),
    @("safe", 0.7605633802816901, 108, 108, 1, 0, $false, 34),
    @("special", 0.75, 27, 27, 1, 0, $false, 9),
    @("safety", 0.7450980392156863, 38, 38, 1, 0, $false, 13),
    @("support", 0.7358490566037735, 78, 78, 1, 0, $false, 28),
    @("confidence", 0.7222222222222222, 26, 26, 1, 0, $false, 10),
    @("good", 0.7, 112, 112, 1, 0, $false, 48),
    @("better", 0.6825396825396826, 43, 43, 1, 0, $false, 20),
    @("well", 0.648936170212766, 61, 61, 1, 0, $false, 33),
    @("heroes", 0.6382978723404256, 30, 30, 1, 0, $false, 17),
    @("relief", 0.62, 31, 31, 1, 0, $false, 19),
    @("hand", 0.5926892950391645, 227, 227, 1, 0, $false, 156),
    @("join", 0.5813953488372093, 25, 25, 1, 0, $false, 18),
    @("fresh", 0.5208333333333334, 25, 25, 1, 0, $false, 23),
    @("protect", 0.5205479452054794, 38, 38, 1, 0, $false, 35),
    @("help", 0.5152542372881356, 152, 152, 1, 0, $false, 143),
    @("hope", 0.5076923076923077, 33, 33, 1, 0, $false, 32),
    @("like", 0.4852941176470588, 165, 165, 1, 0, $false, 175),
    @("care", 0.4831460674157304, 43, 43, 1, 0, $false, 46),
    @("please", 0.4309623430962343, 103, 103, 1, 0, $false, 136),
    @("sure", 0.421875, 27, 27, 1, 0, $false, 37),
    @("increase", 0.4102564102564102, 32, 32, 1, 0, $false, 46),
    @("share", 0.3571428571428572, 25, 25, 1, 0, $false, 45),
    @("online", 0.1004784688995215, 42, 42, 1, 0, $false, 376),
    @("shopping", 0.07211538461538461, 30, 30, 1, 0, $false, 386),
    @("consumer", 0.05701078582434515, 37, 38, 0.97, 0.03000000000000003, $true, 612),
    @("grocery", 0.05660377358490566, 51, 51, 1, 0, $false, 850),
    @("store", 0.05257270693512305, 47, 47, 1, 0, $false, 847),
    @("supermarket", 0.05074971164936563, 44, 44, 1, 0, $false, 823),
    @("19", 0.04781997187060478, 102, 112, 0.91, 0.08999999999999997, $true, 2031),
    @("co", 0.02845134173941158, 88, 102, 0.86, 0.14, $true, 3005),
    @("corona", 0.0225705329153605, 72, 83, 0.87, 0.13, $true, 3118)
)

# Copy the existing bold/bordered/centered format (style index 1, seen on J3) so that
# any newly-created rows (41, 42) pick up the same look the rest of column J uses.
$ws.Cells.Item(3, 10).Copy()

for ($i = 0; $i -lt $rightData.Length; $i++) {
    $r = 3 + $i
    $row = $rightData[$i]
    $needsFormat = ($r -gt 40)
    $ws.Cells.Item($r, 10).Value = $row[0]   # J - name
    if ($needsFormat) {
        $ws.Cells.Item($r, 10).PasteSpecial(-4122)  # xlPasteFormats
    }
    $ws.Cells.Item($r, 11).Value = $row[1]   # K - anchor score
    $ws.Cells.Item($r, 12).Value = $row[2]   # L - type occurences
    $ws.Cells.Item($r, 13).Value = $row[3]   # M - total occurences
    $ws.Cells.Item($r, 14).Value = $row[4]   # N - +%
    $ws.Cells.Item($r, 15).Value = $row[5]   # O - -%
    $ws.Cells.Item($r, 16).Value = $row[6]   # P - both
    $ws.Cells.Item($r, 17).Value = $row[7]   # Q - normal
}
